# Fills in the "MANTENIMIENTO PREVENTIVO" maintenance log (rows 35-39) on
# the "Computadores" sheet with the completed entries: date (col A),
# activity description (col C) and technician (col M).
#
# Column A receives date-looking text ("2025-05-29", ...) and C37 receives
# a purely-numeric-looking string ("2222"). Assigning those to .Value
# directly would make Excel "smart type" them into a real date serial /
# number (and, for dates, reformat the cell), which would change both the
# stored value and the cell style away from the original style (s="18").
# To keep them as literal text with the original style untouched we stage
# the literal value in an unused scratch cell (prefixed with a leading
# apostrophe so Excel treats it as text), copy it, and paste-special
# "values only" into the destination -- a values-only paste transplants
# the already-resolved text as-is, without re-running Excel's type
# inference, and without touching the destination cell's style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("ZZ1")

function Set-LiteralText($rangeAddress, $text) {
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $ws.Range($rangeAddress).PasteSpecial(-4163)  # xlPasteValues
}

# Row 35
Set-LiteralText "A35" "2025-05-29"
$ws.Range("C35").Value = "ultimo cambio de pantalla"
$ws.Range("M35").Value = "SV Romero Romero Miguel Ángel"

# Row 36
Set-LiteralText "A36" "2025-05-16"
$ws.Range("C36").Value = "opcion 3"
$ws.Range("M36").Value = "SV Romero Romero Miguel Ángel"

# Row 37
Set-LiteralText "A37" "2025-05-22"
Set-LiteralText "C37" "2222"
$ws.Range("M37").Value = "SV Romero Romero Miguel Ángel"

# Row 38
Set-LiteralText "A38" "2025-05-27"
$ws.Range("C38").Value = "prueva del clon "
$ws.Range("M38").Value = "SV Romero Romero Miguel Ángel"

# Row 39
Set-LiteralText "A39" "2025-05-26"
$ws.Range("C39").Value = "clon 1111111"
$ws.Range("M39").Value = "SV Romero Romero Miguel Ángel"

# Clean up the scratch cell so no stray content/formatting is left behind.
$scratch.Clear()
